# Exogenous values.xlsx — "Fixed some units and added *SOME* descriptions
# for variables" commit.
#
# 1. Add unit/description labels in column C for a handful of rows on the
#    "company" and "hubs" sheets (new shared strings get appended in the
#    order the cells are written: companies, dmnl, dmnl, month, $/hub,
#    $/month).
# 2. Bump "total operation cost of one hub" on the "cost" sheet from 100 to
#    1000.
# 3. Re-point the active selection on every sheet, and leave "emissions" as
#    the active (tab-selected) sheet, matching the saved view state.

$wb = $excel.ActiveWorkbook

# --- company sheet: add unit column (C) ---------------------------------
$wsCompany = $wb.Worksheets.Item("company")
$wsCompany.Range("C1").Value = "companies"
$wsCompany.Range("C2").Value = "dmnl"
$wsCompany.Range("C3").Value = "dmnl"
$wsCompany.Range("C4").Value = "month"

# --- hubs sheet: add unit column (C) for the first two rows -------------
$wsHubs = $wb.Worksheets.Item("hubs")
$wsHubs.Range("C1").Value = "$/hub"
$wsHubs.Range("C2").Value = "$/month"

# --- cost sheet: correct "total operation cost of one hub" --------------
$wsCost = $wb.Worksheets.Item("cost")
$wsCost.Range("B3").Value = 1000

# --- restore each sheet's remembered selection ---------------------------
$wsCompany.Activate()
$wsCompany.Range("D4").Select()

$wsHubs.Activate()
$wsHubs.Range("C3").Select()

$wsDemand = $wb.Worksheets.Item("demand")
$wsDemand.Activate()
$wsDemand.Range("C2").Select()

$wsCost.Activate()
$wsCost.Range("B4").Select()

$wsEmissions = $wb.Worksheets.Item("emissions")
$wsEmissions.Activate()
$wsEmissions.Range("D6").Select()
